# Update the "Yazd Ammunition Manufacturing and Metallurgy Industries" /
# "Yazd" / "Iran" entry (row 6) to "Rabita Trust" / "Pakistan", dropping the
# now-unused "Yazd" address cell (D6) entirely, and move the active
# selection to the newly edited cell E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A6").Value = "Rabita Trust"
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = "Pakistan"

# Reproduce the new selection/view state (was G4, now E6).
$ws.Range("E6").Select()
